$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add the new "3h" entries to column B for rows 7-9 (list -> dictionary change)
$ws.Range("B7").Value = "3h"
$ws.Range("B8").Value = "3h"
$ws.Range("B9").Value = "3h"

# Update the active selection to B9
$ws.Range("B9").Select()

# Configure page setup: A4 paper, portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
